# The deck's slide-master theme ("Integral" design, Red-Violet colour
# scheme, backed by ppt/theme/theme1.xml) is recoloured to the stock
# "Office Theme" palette. PowerPoint surfaces a theme's 12 colour-scheme
# slots - dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink, in that fixed
# order - through Master.Theme.ThemeColorScheme, so each slot is updated
# to the corresponding "Office" RGB value.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# index -> (R, G, B) for the standard "Office" theme palette
$officeColors = @(
    @(0,   0,   0),    # 1  dk1      000000
    @(255, 255, 255),  # 2  lt1      FFFFFF
    @(68,  84,  106),  # 3  dk2      44546A
    @(231, 230, 230),  # 4  lt2      E7E6E6
    @(91,  155, 213),  # 5  accent1  5B9BD5
    @(237, 125, 49),   # 6  accent2  ED7D31
    @(165, 165, 165),  # 7  accent3  A5A5A5
    @(255, 192, 0),    # 8  accent4  FFC000
    @(68,  114, 196),  # 9  accent5  4472C4
    @(112, 173, 71),   # 10 accent6  70AD47
    @(5,   99,  193),  # 11 hlink    0563C1
    @(149, 79,  114)   # 12 folHlink 954F72
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $rgb = $officeColors[$i - 1]
    $r = $rgb[0]
    $g = $rgb[1]
    $b = $rgb[2]
    $colorScheme.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
